$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# Row 20
# ---------------------------------------------------------------
$ws.Range("A20").Value = 43753
$ws.Range("B20").Value = 0.79166666666666663
$ws.Range("C20").Value = 0.95833333333333337
$ws.Range("D20").Value = 30
$ws.Range("E20").Value = 210
$ws.Range("F20").Value = "프로토타입 express framework 적용 / refactoring (html 분리)"

# ---------------------------------------------------------------
# Row 21  (font/style already matches target - 돋움 - nothing to change there)
# ---------------------------------------------------------------
$ws.Range("A21").Value = 43754
$ws.Range("B21").Value = 0.79166666666666663
$ws.Range("C21").Value = 0.875
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 120
$ws.Range("F21").Value = "프로토타입 Code refactoring - DB 연결"

# ---------------------------------------------------------------
# Row 22
# ---------------------------------------------------------------
$ws.Range("A22").Value = 43758
$ws.Range("B22").Value = 0.70833333333333337
$ws.Range("C22").Value = 0.79166666666666663
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 120
$ws.Range("F22").Value = "프로토타입 Code refactoring - DB 연결 및 입력한 DB 저장"

# ---------------------------------------------------------------
# Row 23
# ---------------------------------------------------------------
$ws.Range("A23").Value = 43761
$ws.Range("B23").Value = 0.91666666666666663
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 120
$ws.Range("F23").Value = "프로토타입 Code refactoring - DB 연결 및 입력한 DB 저장"
$ws.Rows.Item(23).RowHeight = 13

# F20, F22 and F23 use the same "Dotum" cell style that F21 already has,
# so copy the cell formatting from F21 across instead of re-deriving a font.
$ws.Range("F21").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------
$ws.Range("A24").Value = 43765
$ws.Range("B24").Value = 0.83333333333333337
$ws.Range("C24").Value = 0.95833333333333337
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 180
$ws.Range("F24").Value = "Initial Data Set 수정 및 최종 SRS 1.0 제출"
$len1 = ("Initial Data Set 수정").Length
$total = ("Initial Data Set 수정 및 최종 SRS 1.0 제출").Length
$run2 = $ws.Range("F24").Characters($len1 + 1, $total - $len1)
$run2.Font.Name = "돋움"
$run2.Font.Size = 10
$run2.Font.ColorIndex = -4105

# ---------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------
$ws.Range("A25").Value = 43766
$ws.Range("B25").Value = 0.83333333333333337
$ws.Range("C25").Value = 0.083333333333333329
$ws.Range("D25").Value = 60
$ws.Range("E25").Value = 300
$ws.Range("F25").Value = "SRS 발표자료 PPT 준비"
$len1 = ("SRS 발표자료").Length
$total = ("SRS 발표자료 PPT 준비").Length
$run2 = $ws.Range("F25").Characters($len1 + 1, $total - $len1)
$run2.Font.Name = "돋움"
$run2.Font.Size = 10
$run2.Font.ColorIndex = -4105
$ws.Rows.Item(25).RowHeight = 13

# ---------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------
$ws.Range("A26").Value = 43767
$ws.Range("B26").Value = 0.33333333333333331
$ws.Range("C26").Value = 0.5
$ws.Range("D26").Value = 30
$ws.Range("E26").Value = 210
$ws.Range("F26").Value = "SRS 발표자료 보완 및 연습"
$p0 = ("SRS ").Length
$p1 = ("SRS 발표자료").Length
$p2 = ("SRS 발표자료 보완").Length
$p3 = ("SRS 발표자료 보완 및 연습").Length
$run2 = $ws.Range("F26").Characters($p0 + 1, $p1 - $p0)
$run2.Font.Name = "돋움"
$run2.Font.Size = 10
$run2.Font.ColorIndex = -4105
$run3 = $ws.Range("F26").Characters($p1 + 1, $p2 - $p1)
$run3.Font.Name = "Arial"
$run3.Font.Size = 10
$run3.Font.ColorIndex = -4105
$run4 = $ws.Range("F26").Characters($p2 + 1, $p3 - $p2)
$run4.Font.Name = "돋움"
$run4.Font.Size = 10
$run4.Font.ColorIndex = -4105
$ws.Rows.Item(26).RowHeight = 13

# ---------------------------------------------------------------
# Sheet view: scroll position + selection
# ---------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A29").Select()
